$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44875
$ws.Cells.Item(2, 11).Value = 'Castle Brite'
$ws.Cells.Item(2, 12).Value = 'Primera'
$ws.Cells.Item(2, 13).Value = 50
$ws.Cells.Item(2, 14).Value = 31000
$ws.Cells.Item(2, 15).Value = 32000
$ws.Cells.Item(2, 16).Value = 31400
$ws.Cells.Item(2, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(2, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(2, 19).Value = 3140
$ws.Cells.Item(2, 20).Value = 10

# Row 3
$ws.Cells.Item(3, 4).Value = 44944
$ws.Cells.Item(3, 11).Value = 'Modesto'
$ws.Cells.Item(3, 12).Value = 'Primera'
$ws.Cells.Item(3, 13).Value = 100
$ws.Cells.Item(3, 14).Value = 16000
$ws.Cells.Item(3, 15).Value = 17000
$ws.Cells.Item(3, 16).Value = 16500
$ws.Cells.Item(3, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(3, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(3, 19).Value = 1031
$ws.Cells.Item(3, 20).Value = 16

# Row 4
$ws.Cells.Item(4, 4).Value = 44944
$ws.Cells.Item(4, 11).Value = 'Modesto'
$ws.Cells.Item(4, 12).Value = 'Segunda'
$ws.Cells.Item(4, 13).Value = 50
$ws.Cells.Item(4, 14).Value = 14000
$ws.Cells.Item(4, 15).Value = 14000
$ws.Cells.Item(4, 16).Value = 14000
$ws.Cells.Item(4, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(4, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(4, 19).Value = 875
$ws.Cells.Item(4, 20).Value = 16

# Row 5
$ws.Cells.Item(5, 4).Value = 44902
$ws.Cells.Item(5, 11).Value = 'Castle Brite'
$ws.Cells.Item(5, 12).Value = 'Primera'
$ws.Cells.Item(5, 13).Value = 100
$ws.Cells.Item(5, 14).Value = 15000
$ws.Cells.Item(5, 15).Value = 16000
$ws.Cells.Item(5, 16).Value = 15500
$ws.Cells.Item(5, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(5, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(5, 19).Value = 1550
$ws.Cells.Item(5, 20).Value = 10

# Row 6
$ws.Cells.Item(6, 4).Value = 44902
$ws.Cells.Item(6, 11).Value = 'Castle Brite'
$ws.Cells.Item(6, 12).Value = 'Segunda'
$ws.Cells.Item(6, 13).Value = 50
$ws.Cells.Item(6, 14).Value = 13000
$ws.Cells.Item(6, 15).Value = 13000
$ws.Cells.Item(6, 16).Value = 13000
$ws.Cells.Item(6, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(6, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(6, 19).Value = 1300
$ws.Cells.Item(6, 20).Value = 10

# Row 7
$ws.Cells.Item(7, 4).Value = 44159
$ws.Cells.Item(7, 11).Value = 'Castle Brite'
$ws.Cells.Item(7, 12).Value = 'Primera'
$ws.Cells.Item(7, 13).Value = 100
$ws.Cells.Item(7, 14).Value = 14000
$ws.Cells.Item(7, 15).Value = 15000
$ws.Cells.Item(7, 16).Value = 14500
$ws.Cells.Item(7, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(7, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(7, 19).Value = 967
$ws.Cells.Item(7, 20).Value = 15

# Row 8
$ws.Cells.Item(8, 4).Value = 44189
$ws.Cells.Item(8, 11).Value = 'Dina'
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 200
$ws.Cells.Item(8, 14).Value = 15000
$ws.Cells.Item(8, 15).Value = 16000
$ws.Cells.Item(8, 16).Value = 15500
$ws.Cells.Item(8, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(8, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(8, 19).Value = 1033
$ws.Cells.Item(8, 20).Value = 15

# Row 9
$ws.Cells.Item(9, 4).Value = 44189
$ws.Cells.Item(9, 11).Value = 'Dina'
$ws.Cells.Item(9, 12).Value = 'Segunda'
$ws.Cells.Item(9, 13).Value = 100
$ws.Cells.Item(9, 14).Value = 14000
$ws.Cells.Item(9, 15).Value = 14000
$ws.Cells.Item(9, 16).Value = 14000
$ws.Cells.Item(9, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(9, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(9, 19).Value = 933
$ws.Cells.Item(9, 20).Value = 15

# Row 10
$ws.Cells.Item(10, 4).Value = 44559
$ws.Cells.Item(10, 11).Value = 'Modesto'
$ws.Cells.Item(10, 12).Value = 'Primera'
$ws.Cells.Item(10, 13).Value = 100
$ws.Cells.Item(10, 14).Value = 19000
$ws.Cells.Item(10, 15).Value = 20000
$ws.Cells.Item(10, 16).Value = 19500
$ws.Cells.Item(10, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(10, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(10, 19).Value = 1083
$ws.Cells.Item(10, 20).Value = 18

# Row 11
$ws.Cells.Item(11, 4).Value = 44559
$ws.Cells.Item(11, 11).Value = 'Modesto'
$ws.Cells.Item(11, 12).Value = 'Segunda'
$ws.Cells.Item(11, 13).Value = 50
$ws.Cells.Item(11, 14).Value = 18000
$ws.Cells.Item(11, 15).Value = 18000
$ws.Cells.Item(11, 16).Value = 18000
$ws.Cells.Item(11, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(11, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(11, 19).Value = 1000
$ws.Cells.Item(11, 20).Value = 18

# Row 12
$ws.Cells.Item(12, 4).Value = 44545
$ws.Cells.Item(12, 11).Value = 'Castle Brite'
$ws.Cells.Item(12, 12).Value = 'Primera'
$ws.Cells.Item(12, 13).Value = 100
$ws.Cells.Item(12, 14).Value = 18000
$ws.Cells.Item(12, 15).Value = 19000
$ws.Cells.Item(12, 16).Value = 18500
$ws.Cells.Item(12, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(12, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(12, 19).Value = 1233
$ws.Cells.Item(12, 20).Value = 15

# Row 13
$ws.Cells.Item(13, 4).Value = 44545
$ws.Cells.Item(13, 11).Value = 'Castle Brite'
$ws.Cells.Item(13, 12).Value = 'Segunda'
$ws.Cells.Item(13, 13).Value = 50
$ws.Cells.Item(13, 14).Value = 17000
$ws.Cells.Item(13, 15).Value = 17000
$ws.Cells.Item(13, 16).Value = 17000
$ws.Cells.Item(13, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(13, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(13, 19).Value = 1133
$ws.Cells.Item(13, 20).Value = 15

# Row 14
$ws.Cells.Item(14, 4).Value = 44187
$ws.Cells.Item(14, 11).Value = 'Dina'
$ws.Cells.Item(14, 12).Value = 'Primera'
$ws.Cells.Item(14, 13).Value = 100
$ws.Cells.Item(14, 14).Value = 15000
$ws.Cells.Item(14, 15).Value = 16000
$ws.Cells.Item(14, 16).Value = 15500
$ws.Cells.Item(14, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(14, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(14, 19).Value = 861
$ws.Cells.Item(14, 20).Value = 18

# Row 15
$ws.Cells.Item(15, 4).Value = 44908
$ws.Cells.Item(15, 11).Value = 'Albaricoque'
$ws.Cells.Item(15, 12).Value = 'Primera'
$ws.Cells.Item(15, 13).Value = 100
$ws.Cells.Item(15, 14).Value = 20000
$ws.Cells.Item(15, 15).Value = 22000
$ws.Cells.Item(15, 16).Value = 21000
$ws.Cells.Item(15, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(15, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(15, 19).Value = 1167
$ws.Cells.Item(15, 20).Value = 18

# Row 16
$ws.Cells.Item(16, 4).Value = 44918
$ws.Cells.Item(16, 11).Value = 'Dina'
$ws.Cells.Item(16, 12).Value = 'Primera'
$ws.Cells.Item(16, 13).Value = 100
$ws.Cells.Item(16, 14).Value = 17000
$ws.Cells.Item(16, 15).Value = 18000
$ws.Cells.Item(16, 16).Value = 17500
$ws.Cells.Item(16, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(16, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(16, 19).Value = 972
$ws.Cells.Item(16, 20).Value = 18

# Row 17
$ws.Cells.Item(17, 4).Value = 44901
$ws.Cells.Item(17, 11).Value = 'Castle Brite'
$ws.Cells.Item(17, 12).Value = 'Primera'
$ws.Cells.Item(17, 13).Value = 100
$ws.Cells.Item(17, 14).Value = 15000
$ws.Cells.Item(17, 15).Value = 16000
$ws.Cells.Item(17, 16).Value = 15500
$ws.Cells.Item(17, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(17, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(17, 19).Value = 1550
$ws.Cells.Item(17, 20).Value = 10

# Row 18
$ws.Cells.Item(18, 4).Value = 44938
$ws.Cells.Item(18, 11).Value = 'Modesto'
$ws.Cells.Item(18, 12).Value = 'Primera'
$ws.Cells.Item(18, 13).Value = 270
$ws.Cells.Item(18, 14).Value = 14000
$ws.Cells.Item(18, 15).Value = 15000
$ws.Cells.Item(18, 16).Value = 14556
$ws.Cells.Item(18, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(18, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(18, 19).Value = 970
$ws.Cells.Item(18, 20).Value = 15

# Row 19
$ws.Cells.Item(19, 4).Value = 44579
$ws.Cells.Item(19, 11).Value = 'Modesto'
$ws.Cells.Item(19, 12).Value = 'Primera'
$ws.Cells.Item(19, 13).Value = 180
$ws.Cells.Item(19, 14).Value = 13000
$ws.Cells.Item(19, 15).Value = 14000
$ws.Cells.Item(19, 16).Value = 13444
$ws.Cells.Item(19, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(19, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(19, 19).Value = 747
$ws.Cells.Item(19, 20).Value = 18
